# Refactor synthetic array /3
# Update status icons/labels:
#   - A2: "green square" emoji -> "green book" emoji (statut stays "vert")
#   - A3/A4: "black square" emoji -> "blue book" emoji
#   - B3/B4: "noir" -> "bleu"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "📗"

$ws.Range("A3").Value = "📘"
$ws.Range("A4").Value = "📘"

$ws.Range("B3").Value = "bleu"
$ws.Range("B4").Value = "bleu"
